$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.27

# Row 5
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 0.05

# Row 6
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 52

# Row 7
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0.18

# Row 11
$ws.Range("C11").Value = 44

# Row 13
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 0.18

# Row 14
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0.04

# Row 15
$ws.Range("C15").Value = 8

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0.37

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = 0.37

# Row 18
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 1

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0.06

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.37

# Row 21
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 0.09

# Row 25
$ws.Range("C25").Value = 0
$ws.Range("E25").Value = 1

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0.15

# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0.18

# Row 30
$ws.Range("C30").Value = 0
$ws.Range("E30").Value = 1

# Row 35
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 0.16

# Row 36
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 0.01

# Row 37
$ws.Range("C37").Value = 7
$ws.Range("D37").Value = 10
$ws.Range("E37").Value = 0.07000000000000001
